# "add more yellow data" -- append newly-measured current_y / unc_y (columns
# H & I) readings for several existing voltage rows, bump the F25 (current_b)
# reading, and leave the selection where the user last clicked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing "blue" column point that was re-measured ---------------------
$ws.Range("F25").Value = 23432

# --- new "yellow" (current_y / unc_y) readings ------------------------------
# Row 12 (voltage -2): scientific-format cells, matching the rest of column H/I.
$ws.Range("H12").Value = -1.361
$ws.Range("H12").NumberFormat = "0.00E+00"
$ws.Range("I12").Value = 0.115
$ws.Range("I12").NumberFormat = "0.00E+00"

# Row 20 (voltage -0.9): entered as plain/General numbers (no special format).
$ws.Range("H20").Value = 6.7190000000000003
$ws.Range("I20").Value = 0.0766

# Row 21 (voltage -0.8): H21 keeps the scientific format, I21 stays General.
$ws.Range("H21").Value = 59.1
$ws.Range("H21").NumberFormat = "0.00E+00"
$ws.Range("I21").Value = 1.9059999999999999

# Row 22 (voltage -0.75)
$ws.Range("H22").Value = 140.9
$ws.Range("H22").NumberFormat = "0.00E+00"
$ws.Range("I22").Value = 0.94399999999999995
$ws.Range("I22").NumberFormat = "0.00E+00"

# Row 23 (voltage -0.5)
$ws.Range("H23").Value = 1403
$ws.Range("H23").NumberFormat = "0.00E+00"
$ws.Range("I23").Value = 57.9
$ws.Range("I23").NumberFormat = "0.00E+00"

# Row 24 (voltage -0.25)
$ws.Range("H24").Value = 3494
$ws.Range("H24").NumberFormat = "0.00E+00"
$ws.Range("I24").Value = 196
$ws.Range("I24").NumberFormat = "0.00E+00"

# --- leave the cursor where the author last left it -------------------------
$ws.Range("M26").Select()
